$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '39.288.20'
$ws.Range('E2').Value = '  +1.53%  '

$ws.Range('D3').Value = '2.152.82'
$ws.Range('E3').Value = '  +3.23%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.18%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '229.09'
$ws.Range('E5').Value = '  +0.35%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.621'
$ws.Range('E6').Value = '  +1.65%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '62.67'
$ws.Range('E7').Value = '  +3.57%  '

$ws.Range('E8').Value = '  -0.04%  '

$ws.Range('E9').Value = '  +2.27%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0853'
$ws.Range('E10').Value = '  +2.02%  '

$ws.Range('E11').Value = '  -0.34%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '15.99'
$ws.Range('E12').Value = '  +7.01%  '

$ws.Range('D13').Value = '2.470.44'
$ws.Range('E13').Value = '  +3.09%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '22.29'
$ws.Range('E14').Value = '  +1.88%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.815'
$ws.Range('E15').Value = '  +2.54%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '5.55'
$ws.Range('E16').Value = '  +1.33%  '

$ws.Range('D17').Value = '2.135.47'
$ws.Range('E17').Value = '  +2.26%  '

$ws.Range('D18').Value = '39.312.81'
$ws.Range('E18').Value = '  +1.65%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '72.08'
$ws.Range('E19').Value = '  +0.67%  '

$ws.Range('E20').Value = '  +2.09%  '

$ws.Range('D21').Value = '0.0₃0853'
$ws.Range('E21').Value = '  +1.79%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '228.23'
$ws.Range('E22').Value = '  +0.81%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  -0.10%  '

$ws.Range('E24').Value = '  +1.18%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.34'
$ws.Range('E25').Value = '  -0.02%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '9.74'
$ws.Range('E26').Value = '  +3.26%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '170.65'
$ws.Range('E27').Value = '  -0.01%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.138'
$ws.Range('E28').Value = '  +0.19%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '19.66'
$ws.Range('E29').Value = '  +2.63%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '2.57'

$ws.Range('E32').Value = '  +0.88%  '

$ws.Range('E33').Value = '  +2.19%  '

$ws.Range('E34').Value = '  +2.45%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '7.15'
$ws.Range('E35').Value = '  +11.85%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.0619'
$ws.Range('E36').Value = '  +0.90%  '

$ws.Range('E37').Value = '  +0.94%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.56'
$ws.Range('E38').Value = '  +0.71%  '

$ws.Range('E39').Value = '  -0.08%  '

$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.0230'
$ws.Range('E40').Value = '  +3.05%  '

$ws.Range('B41').Value = 'InjectiveProtocol'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '18.01'
$ws.Range('E41').Value = '  -1.22%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '103.04'
$ws.Range('E42').Value = '  +2.25%  '

$ws.Range('D43').Value = '1.535.75'
$ws.Range('E43').Value = '  -0.23%  '

$ws.Range('E44').Value = '  +6.36%  '

$ws.Range('B45').Value = 'ARBITRUM'
$ws.Range('C45').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.10'
$ws.Range('E45').Value = '  +7.04%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '7.85'
$ws.Range('E46').Value = '  +2.04%  '

$ws.Range('B47').Value = 'HuobiToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.81'
$ws.Range('E47').Value = '  -0.30%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0919'
$ws.Range('E48').Value = '  -0.49%  '

$ws.Range('E49').Value = '  +1.10%  '

$ws.Range('B50').Value = 'MXToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.97'
$ws.Range('E50').Value = '  +0.21%  '

$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').Value = '2.353.40'
$ws.Range('E51').Value = '  +2.94%  '

